$d = $word.ActiveDocument

# The hyperlink text "Automate the Boring Stuff" (in the "Working with Data
# in Python" workshop row) is split across three runs ("Au", "t",
# "omate the Boring Stuff") left over from earlier edits. Re-typing the
# whole phrase over itself collapses it back into a single run while
# keeping the existing Hyperlink character formatting.
$d.Content.Find.Execute("Automate the Boring Stuff", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Automate the Boring Stuff", 2)
